$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F ("想去人数")
$updates = @{
    2  = 11487
    3  = 10947
    5  = 5
    6  = 990
    8  = 57
    9  = 38
    11 = 10617
    12 = 4103
    14 = 2454
    16 = 34
    17 = 110
    19 = 11100
    20 = 10860
    23 = 8
    25 = 20
}

# Apply the same updates to both the "展览" and "全部类型" sheets
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
